# Auto-generated Excel COM-interop edit script
# Mirrors the XML diff: delete 4 summary rows (33-36), then update
# line-item data (rows 8-32) and the recomputed Grand Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edit: delete the old duplicate summary rows 33:36 ---
# (old rows 37:40 shift up to become the new 33:36)
$ws.Rows("33:36").Delete()

# --- Cell value edits ---

$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = 82
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3'
$ws.Range("E8").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F8").Value = 472
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '38704.00'
$ws.Range("C9").Value = 47
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '4'
$ws.Range("E9").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F9").Value = 662
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '31114.00'
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = ''
$ws.Range("C10").Value = 63
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.0'
$ws.Range("E10").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '0.00'
$ws.Range("A11").Value = 'Each'
$ws.Range("C11").Value = 3
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.0'
$ws.Range("E11").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 50
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '150.00'
$ws.Range("A12").Value = 'Each'
$ws.Range("C12").Value = 22
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.0'
$ws.Range("E12").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 78
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '1716.00'
$ws.Range("C13").Value = 100
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.0'
$ws.Range("E13").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 30
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '3000.00'
$ws.Range("C14").Value = 41
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.0'
$ws.Range("E14").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 219
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '8979.00'
$ws.Range("C15").Value = 69
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.0'
$ws.Range("E15").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F15").Value = 303
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '20907.00'
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = ''
$ws.Range("C16").Value = 31
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.0'
$ws.Range("E16").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '0.00'
$ws.Range("A17").Value = 'R. mtr.'
$ws.Range("C17").Value = 16
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16'
$ws.Range("E17").Value = '20 mm'
$ws.Range("F17").Value = 40
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '640.00'
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = ''
$ws.Range("C18").Value = 62
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17'
$ws.Range("E18").Value = '25 mm'
$ws.Range("F18").Value = 56
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '3472.00'
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = ''
$ws.Range("C19").Value = 54
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.0'
$ws.Range("E19").Value = 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F19").Value = 0
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '0.00'
$ws.Range("A20").Value = 'Mtr.'
$ws.Range("C20").Value = 85
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19'
$ws.Range("E20").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F20").Value = 81
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '6885.00'
$ws.Range("A21").Value = 'Set'
$ws.Range("C21").Value = 35
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.0'
$ws.Range("E21").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F21").Value = 5733
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '200655.00'
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = ''
$ws.Range("C22").Value = 24
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.0'
$ws.Range("E22").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F22").Value = 0
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '0.00'
$ws.Range("C23").Value = 93
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.0'
$ws.Range("E23").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("A24").Value = 'Each'
$ws.Range("C24").Value = 70
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27'
$ws.Range("E24").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F24").Value = 492
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '34440.00'
$ws.Range("C25").Value = 9
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '17.0'
$ws.Range("E25").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = ''
$ws.Range("C26").Value = 12
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '31'
$ws.Range("E26").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F26").Value = 0
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '0.00'
$ws.Range("A27").Value = 'Each'
$ws.Range("C27").Value = 37
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '32'
$ws.Range("E27").Value = ' 50/63 A rating'
$ws.Range("F27").Value = 900
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '33300.00'
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = ''
$ws.Range("C28").Value = 93
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34'
$ws.Range("E28").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range("F28").Value = 0
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '0.00'
$ws.Range("A29").Value = 'Each'
$ws.Range("C29").Value = 22
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35'
$ws.Range("E29").Value = '8 Way (8+2)'
$ws.Range("F29").Value = 2184
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '48048.00'
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = ''
$ws.Range("C30").Value = 66
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36'
$ws.Range("E30").Value = 'Total'
$ws.Range("F30").Value = 0
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '0.00'
$ws.Range("A31").Value = '%'
$ws.Range("C31").Value = 39
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37'
$ws.Range("E31").Value = 'Add Tender Premium '
$ws.Range("C32").Value = 35
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38'
$ws.Range("E32").Value = 'Grand Total'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '432010.00'
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '432010.00'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '432010.00'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '432010.00'
